$d = $word.ActiveDocument

# Locate the paragraph that ends with "LOQ4037: Química Orgânica I (Requisito fraco)"
# (the anchor we keep) and the paragraph containing the trailing
# "Creative Commons Attribution" copyright/footer text (the last paragraph we remove),
# by scanning paragraph text rather than hard-coding indices.
$n = $d.Paragraphs.Count

$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $n; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*LOQ4037*Requisito fraco*") {
        $startIdx = $i
    }
    if ($t -like "*Creative Commons Attribution*") {
        $endIdx = $i
    }
}

if ($startIdx -gt 0 -and $endIdx -ge $startIdx) {
    $pStart = $d.Paragraphs.Item($startIdx)
    $pEnd = $d.Paragraphs.Item($endIdx)

    # Range spanning from just after the "LOQ4037..." paragraph's own
    # paragraph mark through the end of the "© 2020 ... Creative Commons
    # Attribution" paragraph (inclusive of its paragraph mark). This removes
    # the blank paragraph, the "Ver no Jupiter..." paragraph, and the
    # copyright/footer paragraph in one shot, leaving "LOQ4037..." followed
    # directly by whatever paragraph originally came after the footer.
    $r = $d.Range($pStart.Range.End, $pEnd.Range.End)
    $r.Delete()
}
